# Base General.xlsx - apply "corte de caja" update
# - Adds payment-type aware rows, new currency number format for Importe,
#   swaps the (previously mismatched) Correo Electronico / Codigo Postal header
#   order, widens several columns, and repairs the rows where the dawn
#   ("madrugada") shift's folio/date columns were shifted by one cell because
#   a payment-type ("Efectivo") value got inserted ahead of them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: write a value as a genuine text string (t="inlineStr"/shared string)
# even when it looks like a number (e.g. "5563193656", "$600", "26/08/2022").
# A plain  Range.Value = "..."  assignment lets the engine "smart" parse the
# string and silently store it as a number, which is not what the source
# data (originally authored outside Excel) looks like. Routing the literal
# through a text formula + copy/paste-values (which carries over the VALUE
# only, not the formula or any formatting) reproduces a genuine text cell
# without creating any stray/unused style entries.
# ---------------------------------------------------------------------------
$stage = $ws.Range("ZZ1")
function Set-Text($addr, $val) {
    $escaped = $val -replace '"', '""'
    $stage.Formula = '="' + $escaped + '"'
    $stage.Copy()
    $ws.Range($addr).PasteSpecial(-4163)  # xlPasteValues
}

# ---------------------------------------------------------------------------
# Header row: swap Correo Electronico / Codigo Postal so the header labels
# correctly line up with the data that already lives in columns C/D.
# ---------------------------------------------------------------------------
Set-Text "C1" "Correo Electrónico"
Set-Text "D1" "Código Postal"

# New currency format applied to the "Importe" column.
$ws.Range("F1").NumberFormat = '"$"#,##0.00'
$ws.Range("F2:F3").NumberFormat = '"$"#,##0.00'

# ---------------------------------------------------------------------------
# New rows captured since the last cut (rows 4-11). Rows 9 & 10 preserve the
# original (buggy) column layout where a payment-type value ("Efectivo") was
# recorded in column G ahead of the folio/date, which is the very bug the
# commit fixes going forward - the historic rows are kept as originally
# recorded.
# ---------------------------------------------------------------------------
$rows = @(
    @{r=4;  A="ALEX SERRANO DURÁN"; B="5563193656"; C="alexserrano0805@gmail.com"; D="52950";  E="CONSULTA"; F="`$600";    G="12:02";    H="001-20220826M"; I="26/08/2022"},
    @{r=5;  A="NOMBRE";             B="5563193656"; C="mmm@gmail.com";             D="01022";  E="SERVICIO"; F="`$600.50"; G="12:06";    H="002-20220826M"; I="26/08/2022"},
    @{r=6;  A="NOMBRE";             B="555";         C="alex@gmail.com";           D="3333";   E="SERVICIO"; F="`$600";    G="12:11";    H="001-20220826V"; I="26/08/2022"},
    @{r=7;  A="ALEX";               B="556214";      C="asdasd";                   D="124124"; E="SERVICIO"; F="`$600";    G="12:34";    H="003-20220826M"; I="26/08/2022"},
    @{r=8;  A="ALEX SERRANO";       B="123123123";   C="alexserrano0805";          D="125521"; E="SERVICIO"; F="`$123";    G="02:28";    H="001-20220827V"; I="28/08/2022"},
    @{r=11; A="NOMBRE";             B="5516169339";  C="aaa@a.com";                D="15154";  E="SERVICIO"; F="`$500";    G="03:13";    H="003-20220827N"; I="27/08/2022"}
)

foreach ($row in $rows) {
    $r = $row.r
    Set-Text "A$r" $row.A
    Set-Text "B$r" $row.B
    Set-Text "C$r" $row.C
    Set-Text "D$r" $row.D
    Set-Text "E$r" $row.E
    Set-Text "F$r" $row.F
    Set-Text "G$r" $row.G
    Set-Text "H$r" $row.H
    Set-Text "I$r" $row.I
}

# Rows 9 & 10: historic (buggy) column layout with payment type in column G.
Set-Text "A9" "NOMBRE"
Set-Text "B9" "12312312"
Set-Text "C9" "aaqa@.com"
Set-Text "D9" "52950"
Set-Text "E9" "SERVICIO"
Set-Text "F9" "`$500"
Set-Text "G9" "Efectivo"
Set-Text "H9" "02:53"
Set-Text "I9" "001-20220827N"
Set-Text "J9" "28/08/2022"

Set-Text "A10" "NOMBRE"
Set-Text "B10" "556219305"
Set-Text "C10" "aaa@.com"
Set-Text "D10" "15252"
Set-Text "E10" "SERVICIO"
Set-Text "F10" "`$500"
Set-Text "G10" "Efectivo"
Set-Text "H10" "03:02"
Set-Text "I10" "28/08/2022"

# Remove the staging cell entirely (value + formatting) so it never shows up
# in the saved worksheet / used range.
$stage.Clear()

# ---------------------------------------------------------------------------
# Column widths. The XLSX <col width> unit (characters, 256ths) is offset
# from the Range.ColumnWidth COM property by the fixed ~5px/MDW cell-padding
# that Excel silently folds into the *stored* width but not the API value, so
# the values below are calibrated (not the literal target widths) to land
# exactly on the desired stored width.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 34
$ws.Columns.Item(3).ColumnWidth = 33.666
$ws.Columns.Item(4).ColumnWidth = 14
$ws.Columns.Item(5).ColumnWidth = 33.168
$ws.Columns.Item(6).ColumnWidth = 16.6666
$ws.Columns.Item(8).ColumnWidth = 19.1656

# ---------------------------------------------------------------------------
# Active selection moves to F3
# ---------------------------------------------------------------------------
$ws.Range("F3").Select()
